$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "63.411.56"
Set-TextValue $ws.Range("E2") "  +5.46%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.064.47"
Set-TextValue $ws.Range("E3") "  +4.44%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.14%  "

# Row 5
Set-TextValue $ws.Range("D5") "550.34"
Set-TextValue $ws.Range("E5") "  +6.46%  "

# Row 6
Set-TextValue $ws.Range("D6") "138.99"
Set-TextValue $ws.Range("E6") "  +8.39%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.22%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.058.91"
Set-TextValue $ws.Range("E8") "  +4.67%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.500"
Set-TextValue $ws.Range("E9") "  +5.25%  "

# Row 10
Set-TextValue $ws.Range("B10") "Toncoin"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D10") "6.24"
Set-TextValue $ws.Range("E10") "  +3.34%  "

# Row 11
Set-TextValue $ws.Range("B11") "Dogecoin"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D11") "0.150"
Set-TextValue $ws.Range("E11") "  +3.12%  "

# Row 12
Set-TextValue $ws.Range("E12") "  +5.31%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000228"
Set-TextValue $ws.Range("E13") "  +5.66%  "

# Row 14
Set-TextValue $ws.Range("D14") "34.91"
Set-TextValue $ws.Range("E14") "  +7.51%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.565.50"
Set-TextValue $ws.Range("E15") "  +4.19%  "

# Row 16
Set-TextValue $ws.Range("D16") "63.492.94"
Set-TextValue $ws.Range("E16") "  +5.34%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.067.51"
Set-TextValue $ws.Range("E17") "  +4.02%  "

# Row 18
Set-TextValue $ws.Range("E18") "  -0.35%  "

# Row 19
Set-TextValue $ws.Range("E19") "  +5.93%  "

# Row 20
Set-TextValue $ws.Range("D20") "483.12"
Set-TextValue $ws.Range("E20") "  +7.24%  "

# Row 21
Set-TextValue $ws.Range("D21") "13.57"
Set-TextValue $ws.Range("E21") "  +5.85%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.683"
Set-TextValue $ws.Range("E22") "  +3.58%  "

# Row 23
Set-TextValue $ws.Range("D23") "7.22"
Set-TextValue $ws.Range("E23") "  +8.23%  "

# Row 24
Set-TextValue $ws.Range("D24") "81.49"
Set-TextValue $ws.Range("E24") "  +5.77%  "

# Row 25
Set-TextValue $ws.Range("D25") "12.60"
Set-TextValue $ws.Range("E25") "  +9.40%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +0.20%  "

# Row 27
Set-TextValue $ws.Range("D27") "2.76"
Set-TextValue $ws.Range("E27") "  +6.75%  "

# Row 28
Set-TextValue $ws.Range("D28") "7.99"
Set-TextValue $ws.Range("E28") "  +6.38%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.00"
Set-TextValue $ws.Range("E29") "  +11.69%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("E30") "  +0.00%  "

# Row 31
Set-TextValue $ws.Range("D31") "26.05"
Set-TextValue $ws.Range("E31") "  +4.93%  "

# Row 32
Set-TextValue $ws.Range("E32") "  +4.14%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.42"
Set-TextValue $ws.Range("E33") "  +8.64%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.74"
Set-TextValue $ws.Range("E34") "  +9.50%  "

# Row 35
Set-TextValue $ws.Range("D35") "55.43"
Set-TextValue $ws.Range("E35") "  +2.54%  "

# Row 36
Set-TextValue $ws.Range("E36") "  +6.14%  "

# Row 37
Set-TextValue $ws.Range("D37") "467.77"
Set-TextValue $ws.Range("E37") "  +4.90%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.159.21"
Set-TextValue $ws.Range("E38") "  +1.19%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0816"
Set-TextValue $ws.Range("E39") "  +7.17%  "

# Row 40
Set-TextValue $ws.Range("E40") "  +5.96%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.120"
Set-TextValue $ws.Range("E41") "  +4.30%  "

# Row 42
Set-TextValue $ws.Range("D42") "8.27"
Set-TextValue $ws.Range("E42") "  +5.60%  "

# Row 43
Set-TextValue $ws.Range("B43") "dogwifhat"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D43") "2.62"
Set-TextValue $ws.Range("E43") "  +9.67%  "

# Row 44
Set-TextValue $ws.Range("B44") "InjectiveProtocol"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D44") "28.54"
Set-TextValue $ws.Range("E44") "  +16.30%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.253"
Set-TextValue $ws.Range("E45") "  +6.30%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.06"
Set-TextValue $ws.Range("E47") "  +8.80%  "

# Row 48
Set-TextValue $ws.Range("E48") "  +4.28%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0₃0511"
Set-TextValue $ws.Range("E49") "  +2.97%  "

# Row 50
Set-TextValue $ws.Range("D50") "116.57"
Set-TextValue $ws.Range("E50") "  -1.00%  "

# Row 51
Set-TextValue $ws.Range("E51") "  +8.63%  "
